$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing (only) sheet to "Armors" and add two more sheets
#    ("Shields", "Weapons") after it, matching the tab order from the diff.
# ---------------------------------------------------------------------------
$wsArmors = $wb.Worksheets.Item(1)
$wsArmors.Name = "Armors"

$wsShields = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsArmors)
$wsShields.Name = "Shields"

$wsWeapons = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsShields)
$wsWeapons.Name = "Weapons"

# ---------------------------------------------------------------------------
# 2. Armors sheet data
# ---------------------------------------------------------------------------
$wsArmors.Cells.Item(1,1).Value = "name"
$wsArmors.Cells.Item(1,2).Value = "special"
$wsArmors.Cells.Item(1,3).Value = "damage"

$wsArmors.Cells.Item(2,1).Value = "Mała Tarcza"
$wsArmors.Cells.Item(2,2).Value = "Defensywna +1"
$wsArmors.Cells.Item(2,3).Value = 1

$wsArmors.Columns.Item(1).ColumnWidth = 10.25
$wsArmors.Columns.Item(2).ColumnWidth = 13.6

# ---------------------------------------------------------------------------
# 3. Shields sheet data
# ---------------------------------------------------------------------------
$wsShields.Cells.Item(1,1).Value = "name"
$wsShields.Cells.Item(1,2).Value = "special"
$wsShields.Cells.Item(1,3).Value = "damage"

$wsShields.Cells.Item(2,1).Value = "Mała Tarcza"
$wsShields.Cells.Item(2,2).Value = "Defensywna +1"
$wsShields.Cells.Item(2,3).Value = 1

$wsShields.Cells.Item(3,1).Value = "Duża Tarcza"
$wsShields.Cells.Item(3,2).Value = "Defensuwa +2, Rozmiar 1"
$wsShields.Cells.Item(3,3).Value = "1k3"

$wsShields.Columns.Item(1).ColumnWidth = 10.25
$wsShields.Columns.Item(2).ColumnWidth = 22.6

# ---------------------------------------------------------------------------
# 4. Weapons sheet data
# ---------------------------------------------------------------------------
$wsWeapons.Cells.Item(1,1).Value = "name"
$wsWeapons.Cells.Item(1,2).Value = "damage"
$wsWeapons.Cells.Item(1,3).Value = "handle"
$wsWeapons.Cells.Item(1,4).Value = "specials"

$wsWeapons.Cells.Item(2,1).Value = "Pałka"
$wsWeapons.Cells.Item(2,2).Value = "1k6"
$wsWeapons.Cells.Item(2,3).Value = "1H"
$wsWeapons.Cells.Item(2,4).Value = "-"

$wsWeapons.Cells.Item(3,1).Value = "Kostur"
$wsWeapons.Cells.Item(3,2).Value = "1k6+1"
$wsWeapons.Cells.Item(3,3).Value = "2H"
$wsWeapons.Cells.Item(3,4).Value = "Finezyjna"

$wsWeapons.Cells.Item(4,1).Value = "Proca"
$wsWeapons.Cells.Item(4,2).Value = "1k3"
$wsWeapons.Cells.Item(4,3).Value = "1H/2H"
$wsWeapons.Cells.Item(4,4).Value = "Używa kamieni, zasięg (średni)"
$wsWeapons.Cells.Item(4,4).Font.Size = 7
$wsWeapons.Cells.Item(4,4).Font.Name = "Metropolis-Regular"
$wsWeapons.Cells.Item(4,4).Font.Color = 0

$wsWeapons.Cells.Item(5,1).Value = "Sztylet"
$wsWeapons.Cells.Item(5,2).Value = "1k3"
$wsWeapons.Cells.Item(5,3).Value = "1H/2H"
$wsWeapons.Cells.Item(5,4).Value = "Finezyjna, miotana, zasięg (bliski)"

$wsWeapons.Cells.Item(6,1).Value = "Miecz"
$wsWeapons.Cells.Item(6,2).Value = "1k6+2"
$wsWeapons.Cells.Item(6,3).Value = "1H"
$wsWeapons.Cells.Item(6,4).Value = "-"

# ---------------------------------------------------------------------------
# 5. Selections / active sheet to roughly match the saved view state
# ---------------------------------------------------------------------------
$wsArmors.Activate() | Out-Null
$wsArmors.Range("A1:C2").Select() | Out-Null

$wsShields.Activate() | Out-Null
$wsShields.Range("I17").Select() | Out-Null

$wsWeapons.Activate() | Out-Null
$wsWeapons.Range("E10").Select() | Out-Null
